# Update the "380 kV" case results (pl_mw.xlsx / Sheet1) with the newly
# computed per-line active power flow values for rows 2-25 (columns B-O,
# skipping the zero-valued E, J, K, L, N columns that are unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B, rows 2-25
$arr_B = New-Object 'object[,]' 24,1
$arr_B[0,0] = 0.6159107092382214
$arr_B[1,0] = 0.5380219922801643
$arr_B[2,0] = 0.4899664346918939
$arr_B[3,0] = 0.4703271999164258
$arr_B[4,0] = 0.4670627850896949
$arr_B[5,0] = 0.4897017979854184
$arr_B[6,0] = 0.5891039301068872
$arr_B[7,0] = 0.7821140693481539
$arr_B[8,0] = 0.9226562708977326
$arr_B[9,0] = 0.986299576708916
$arr_B[10,0] = 1.010356076524033
$arr_B[11,0] = 1.005177062391624
$arr_B[12,0] = 0.9882796101801432
$arr_B[13,0] = 0.9779236606738664
$arr_B[14,0] = 0.9184910122876886
$arr_B[15,0] = 0.8819551995778738
$arr_B[16,0] = 0.8609135923857707
$arr_B[17,0] = 0.8537846636990025
$arr_B[18,0] = 0.8858473270969967
$arr_B[19,0] = 0.9932440081839218
$arr_B[20,0] = 1.063177441634309
$arr_B[21,0] = 1.025876822987982
$arr_B[22,0] = 0.8840878103391674
$arr_B[23,0] = 0.7301146861109942
$ws.Range("B2:B25").Value = $arr_B

# Column C, rows 2-25
$arr_C = New-Object 'object[,]' 24,1
$arr_C[0,0] = 0.1897266710568601
$arr_C[1,0] = 0.1733923610344732
$arr_C[2,0] = 0.1633164876076023
$arr_C[3,0] = 0.1591992833304232
$arr_C[4,0] = 0.1585149627185558
$arr_C[5,0] = 0.1632610062551976
$arr_C[6,0] = 0.1841045541998767
$arr_C[7,0] = 0.2245894104225954
$arr_C[8,0] = 0.2540722671719209
$arr_C[9,0] = 0.2674231308457422
$arr_C[10,0] = 0.2724695120685112
$arr_C[11,0] = 0.2713831055656044
$arr_C[12,0] = 0.2678384890664915
$arr_C[13,0] = 0.2656660823230368
$arr_C[14,0] = 0.2531984837577284
$arr_C[15,0] = 0.245534009234774
$arr_C[16,0] = 0.2411198923839777
$arr_C[17,0] = 0.2396243829298328
$arr_C[18,0] = 0.2463505007970355
$arr_C[19,0] = 0.2688798854475181
$arr_C[20,0] = 0.2835496854574728
$arr_C[21,0] = 0.2757252922541795
$arr_C[22,0] = 0.2459813889321936
$arr_C[23,0] = 0.2136813574055623
$ws.Range("C2:C25").Value = $arr_C

# Column D, rows 2-25
$arr_D = New-Object 'object[,]' 24,1
$arr_D[0,0] = 0.02191339196804876
$arr_D[1,0] = 0.01910969445793853
$arr_D[2,0] = 0.0173815707484053
$arr_D[3,0] = 0.01667573458360039
$arr_D[4,0] = 0.01655843540295621
$arr_D[5,0] = 0.01737205804388964
$arr_D[6,0] = 0.0209480898591039
$arr_D[7,0] = 0.02790572740148889
$arr_D[8,0] = 0.03298145581587164
$arr_D[9,0] = 0.03528219647709818
$arr_D[10,0] = 0.03615218944526077
$arr_D[11,0] = 0.03596487739173426
$arr_D[12,0] = 0.03535379668933558
$arr_D[13,0] = 0.03497932775665902
$arr_D[14,0] = 0.03283092588199565
$arr_D[15,0] = 0.03151079680542068
$arr_D[16,0] = 0.03075072228205045
$arr_D[17,0] = 0.03049324359389516
$arr_D[18,0] = 0.03165140707540104
$arr_D[19,0] = 0.0355333202701047
$arr_D[20,0] = 0.03806307359080563
$arr_D[21,0] = 0.03671358570213101
$arr_D[22,0] = 0.03158784063894871
$arr_D[23,0] = 0.02602964878838776
$ws.Range("D2:D25").Value = $arr_D

# Column F, rows 2-25
$arr_F = New-Object 'object[,]' 24,1
$arr_F[0,0] = 0.3052248113244787
$arr_F[1,0] = 0.3028798491082725
$arr_F[2,0] = 0.3017303346825528
$arr_F[3,0] = 0.3013347176755374
$arr_F[4,0] = 0.3012734186244757
$arr_F[5,0] = 0.3017247046492315
$arr_F[6,0] = 0.3043559067803372
$arr_F[7,0] = 0.3118279927343721
$arr_F[8,0] = 0.3187408480103997
$arr_F[9,0] = 0.3221976521285441
$arr_F[10,0] = 0.3235517474076488
$arr_F[11,0] = 0.323258110749002
$arr_F[12,0] = 0.3223081497939901
$arr_F[13,0] = 0.3217321476784747
$arr_F[14,0] = 0.3185212351075108
$arr_F[15,0] = 0.3166315242611901
$arr_F[16,0] = 0.3155739782690503
$arr_F[17,0] = 0.3152209489220255
$arr_F[18,0] = 0.3168296465356519
$arr_F[19,0] = 0.3225859513670457
$arr_F[20,0] = 0.326610872604455
$arr_F[21,0] = 0.3244385812989137
$arr_F[22,0] = 0.316739985534376
$arr_F[23,0] = 0.3095576503645532
$ws.Range("F2:F25").Value = $arr_F

# Column G, rows 2-25
$arr_G = New-Object 'object[,]' 24,1
$arr_G[0,0] = 0.1662234647975112
$arr_G[1,0] = 0.1649373110469767
$arr_G[2,0] = 0.1643614004117993
$arr_G[3,0] = 0.1641802016376701
$arr_G[4,0] = 0.1641533364143584
$arr_G[5,0] = 0.1643587404842464
$arr_G[6,0] = 0.1657354627813561
$arr_G[7,0] = 0.1701441827593868
$arr_G[8,0] = 0.1744431792544177
$arr_G[9,0] = 0.1766329622975107
$arr_G[10,0] = 0.1774961475698902
$arr_G[11,0] = 0.1773087300727028
$arr_G[12,0] = 0.1767032947529188
$arr_G[13,0] = 0.1763368796327498
$arr_G[14,0] = 0.1743048083984178
$arr_G[15,0] = 0.1731183772859737
$arr_G[16,0] = 0.1724579952249385
$arr_G[17,0] = 0.1722381745544439
$arr_G[18,0] = 0.1732423936217842
$arr_G[19,0] = 0.1768802017430033
$arr_G[20,0] = 0.1794558107782365
$arr_G[21,0] = 0.178062938210104
$arr_G[22,0] = 0.1731862582272967
$arr_G[23,0] = 0.1687666781804822
$ws.Range("G2:G25").Value = $arr_G

# Column H, rows 2-25
$arr_H = New-Object 'object[,]' 24,1
$arr_H[0,0] = 0.3356984809392287
$arr_H[1,0] = 0.3386419939305583
$arr_H[2,0] = 0.3406654445999209
$arr_H[3,0] = 0.3415443152383162
$arr_H[4,0] = 0.341693529114302
$arr_H[5,0] = 0.3406770775704544
$arr_H[6,0] = 0.3366685194313348
$arr_H[7,0] = 0.3305248442795516
$arr_H[8,0] = 0.3270611893598456
$arr_H[9,0] = 0.3257142324450939
$arr_H[10,0] = 0.3252371165877719
$arr_H[11,0] = 0.3253384056566588
$arr_H[12,0] = 0.3256743191525757
$arr_H[13,0] = 0.325884368503516
$arr_H[14,0] = 0.3271538232021243
$arr_H[15,0] = 0.3279912077825031
$arr_H[16,0] = 0.3284943695986087
$arr_H[17,0] = 0.3286684257326016
$arr_H[18,0] = 0.3278998390265428
$arr_H[19,0] = 0.3255747585208866
$arr_H[20,0] = 0.324247256887773
$arr_H[21,0] = 0.3249381736675758
$arr_H[22,0] = 0.3279410791359325
$arr_H[23,0] = 0.3320027013085394
$ws.Range("H2:H25").Value = $arr_H

# Column I, rows 2-25
$arr_I = New-Object 'object[,]' 24,1
$arr_I[0,0] = 0.2283417952682925
$arr_I[1,0] = 0.2335664298988398
$arr_I[2,0] = 0.2370092601481257
$arr_I[3,0] = 0.2384711145919205
$arr_I[4,0] = 0.2387174046267599
$arr_I[5,0] = 0.2370287371071758
$arr_I[6,0] = 0.230094385473361
$arr_I[7,0] = 0.2183681924949603
$arr_I[8,0] = 0.2109054533406614
$arr_I[9,0] = 0.2077631192850955
$arr_I[10,0] = 0.2066097140968299
$arr_I[11,0] = 0.2068564924377636
$arr_I[12,0] = 0.207667494605948
$arr_I[13,0] = 0.2081690204414031
$arr_I[14,0] = 0.2111159119770285
$arr_I[15,0] = 0.212988562298154
$arr_I[16,0] = 0.214089414962805
$arr_I[17,0] = 0.2144662178184671
$arr_I[18,0] = 0.2127867554701073
$arr_I[19,0] = 0.2074282904912224
$arr_I[20,0] = 0.2041392417748895
$arr_I[21,0] = 0.2058751054462107
$arr_I[22,0] = 0.2128779168482211
$arr_I[23,0] = 0.2213388610178697
$ws.Range("I2:I25").Value = $arr_I

# Column M, rows 2-25
$arr_M = New-Object 'object[,]' 24,1
$arr_M[0,0] = 0.8958937206540156
$arr_M[1,0] = 0.797093620861375
$arr_M[2,0] = 0.7369500878308202
$arr_M[3,0] = 0.7125658685623932
$arr_M[4,0] = 0.708524248516639
$arr_M[5,0] = 0.7366207364426316
$arr_M[6,0] = 0.8617158016676143
$arr_M[7,0] = 1.111442036355726
$arr_M[8,0] = 1.29803037629442
$arr_M[9,0] = 1.383685535601131
$arr_M[10,0] = 1.416240008561431
$arr_M[11,0] = 1.409223422581661
$arr_M[12,0] = 1.386361390068828
$arr_M[13,0] = 1.372373414796328
$arr_M[14,0] = 1.292448824537701
$arr_M[15,0] = 1.243621321416953
$arr_M[16,0] = 1.215609552215199
$arr_M[17,0] = 1.206137485857781
$arr_M[18,0] = 1.248811535094447
$arr_M[19,0] = 1.393073242345338
$arr_M[20,0] = 1.48805212044563
$arr_M[21,0] = 1.437294019813862
$arr_M[22,0] = 1.246464853733627
$arr_M[23,0] = 1.043368696418526
$ws.Range("M2:M25").Value = $arr_M

# Column O, rows 2-25
$arr_O = New-Object 'object[,]' 24,1
$arr_O[0,0] = 0.9057982176320536
$arr_O[1,0] = 0.9089946253170069
$arr_O[2,0] = 0.9118546664451799
$arr_O[3,0] = 0.913245144734006
$arr_O[4,0] = 0.9134896014193004
$arr_O[5,0] = 0.9118725088576696
$arr_O[6,0] = 0.9067136664298374
$arr_O[7,0] = 0.903749691253438
$arr_O[8,0] = 0.9059780463411897
$arr_O[9,0] = 0.9079583942231864
$arr_O[10,0] = 0.908848071104785
$arr_O[11,0] = 0.9086502352623995
$arr_O[12,0] = 0.9080287831003773
$arr_O[13,0] = 0.9076663493209622
$arr_O[14,0] = 0.9058681424079253
$arr_O[15,0] = 0.9050131248409485
$arr_O[16,0] = 0.9046122826219687
$arr_O[17,0] = 0.9044921604521932
$arr_O[18,0] = 0.9050947253510486
$arr_O[19,0] = 0.9082075195296113
$arr_O[20,0] = 0.9110569015155079
$arr_O[21,0] = 0.9094613110235059
$arr_O[22,0] = 0.9050575512490866
$arr_O[23,0] = 0.9037810758933347
$ws.Range("O2:O25").Value = $arr_O

